# Updated cryptos list on Thu May 11 11:59:19 UTC 2023 with GitHub Actions
#
# Refreshes the Price (D) and Volume(1h) (E) columns on the active sheet.
# Both columns hold plain text in the workbook (not numbers), so any value
# that Excel would otherwise auto-convert to a number on assignment
# (single-dot decimals like "1.002") is written with a leading apostrophe
# to force text entry, exactly as a user typing into a text-looking cell
# would. Values that already can't parse as a number (e.g. "27.462.47",
# with two dots) or the "  +x.xx%  " volume strings are assigned as-is.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.462.47"
$ws.Range("E2").Value = "  -0.91%  "
$ws.Range("D3").Value = "1.825.06"
$ws.Range("E3").Value = "  -1.51%  "
$ws.Range("D4").Value = "'1.002"
$ws.Range("E4").Value = "  +0.02%  "
$ws.Range("E5").Value = "  +0.13%  "
$ws.Range("E6").Value = "  +0.04%  "
$ws.Range("D7").Value = "'0.4260"
$ws.Range("E7").Value = "  -0.20%  "
$ws.Range("E8").Value = "  +0.40%  "
$ws.Range("D9").Value = "'0.07212"
$ws.Range("E9").Value = "  -1.15%  "
$ws.Range("D10").Value = "'0.8601"
$ws.Range("E10").Value = "  -2.15%  "
$ws.Range("D11").Value = "'20.56"
$ws.Range("E11").Value = "  -1.09%  "
$ws.Range("D12").Value = "1.911.05"
$ws.Range("E12").Value = "  +3.01%  "
$ws.Range("E13").Value = "  +0.71%  "
$ws.Range("D14").Value = "'6.470"
$ws.Range("E14").Value = "  -1.22%  "
$ws.Range("D15").Value = "'0.06928"
$ws.Range("E15").Value = "  -0.86%  "
$ws.Range("D16").Value = "'1.004"
$ws.Range("E16").Value = "  +0.01%  "
$ws.Range("D17").Value = "'80.76"
$ws.Range("E17").Value = "  +1.24%  "
$ws.Range("D18").Value = "'0.000008895"
$ws.Range("E18").Value = "  -0.62%  "
$ws.Range("E19").Value = "  +0.04%  "
$ws.Range("D20").Value = "'15.37"
$ws.Range("E20").Value = "  +0.79%  "
$ws.Range("D21").Value = "27.567.16"
$ws.Range("E21").Value = "  -0.26%  "
$ws.Range("D22").Value = "'5.126"
$ws.Range("E22").Value = "  +2.49%  "
$ws.Range("E23").Value = "  +4.39%  "
$ws.Range("D24").Value = "2.099.27"
$ws.Range("E24").Value = "  +1.41%  "
$ws.Range("D25").Value = "'1.988"
$ws.Range("E25").Value = "  +0.00%  "
$ws.Range("D26").Value = "'155.35"
$ws.Range("E26").Value = "  +0.51%  "
$ws.Range("D27").Value = "'18.71"
$ws.Range("E27").Value = "  +1.27%  "
$ws.Range("D28").Value = "'5.131"
$ws.Range("E28").Value = "  -2.29%  "
$ws.Range("D29").Value = "'114.16"
$ws.Range("E29").Value = "  -4.97%  "
$ws.Range("D30").Value = "'1.788"
$ws.Range("E30").Value = "  -4.95%  "
$ws.Range("D31").Value = "'0.08893"
$ws.Range("E31").Value = "  -0.06%  "
$ws.Range("D32").Value = "'2.988"
$ws.Range("E32").Value = "  +0.80%  "
$ws.Range("D33").Value = "'0.7444"
$ws.Range("E33").Value = "  -2.09%  "
$ws.Range("D34").Value = "'4.543"
$ws.Range("E34").Value = "  +0.60%  "
$ws.Range("D35").Value = "'1.117"
$ws.Range("E35").Value = "  -0.40%  "
$ws.Range("E36").Value = "  +0.06%  "
$ws.Range("D37").Value = "'1.085"
$ws.Range("E37").Value = "  -2.14%  "
$ws.Range("D38").Value = "'0.05253"
$ws.Range("E38").Value = "  -3.18%  "
$ws.Range("E39").Value = "  -0.51%  "
$ws.Range("D40").Value = "'2.788"
$ws.Range("E40").Value = "  -1.26%  "
$ws.Range("D41").Value = "'0.5070"
$ws.Range("E41").Value = "  -0.36%  "
$ws.Range("D42").Value = "'0.1651"
$ws.Range("E42").Value = "  -1.27%  "
$ws.Range("D43").Value = "'6.366"
$ws.Range("E43").Value = "  -3.87%  "
$ws.Range("D44").Value = "'8.330"
$ws.Range("E44").Value = "  -1.02%  "
$ws.Range("D45").Value = "'10.43"
$ws.Range("E45").Value = "  +0.49%  "
$ws.Range("D46").Value = "'106.41"
$ws.Range("E46").Value = "  +0.56%  "
$ws.Range("D47").Value = "'0.06454"
$ws.Range("E47").Value = "  -1.17%  "
$ws.Range("D48").Value = "'0.4678"
$ws.Range("E48").Value = "  +0.17%  "
$ws.Range("E49").Value = "  +0.06%  "
$ws.Range("E50").Value = "  -0.54%  "
$ws.Range("D51").Value = "'63.83"
$ws.Range("E51").Value = "  -1.12%  "
